$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($row, $col, $refRow, [string]$text)
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $ws.Cells.Item($refRow, $col).Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null
}

function Set-NumberCell {
    param($row, $col, $refRow, $num)
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = $num
    $ws.Cells.Item($refRow, $col).Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null
}

# ---- Row 31 ----
Set-TextCell 31 1 29 "2/13/2020"
Set-NumberCell 31 2 29 0.993055555555556
Set-TextCell 31 3 29 "None"
Set-TextCell 31 4 29 "Review the slides"
Set-TextCell 31 5 29 "Done"
Set-TextCell 31 6 29 "Today I think I have learnt something outside the codes. It makes me realize that people related to the project, such as stakeholder, developers, are also important for a project"
Set-TextCell 31 7 29 "Not bad"
$ws.Rows.Item(31).RowHeight = 109

# ---- Row 32 ----
Set-TextCell 32 1 29 "2/16/2020"
Set-NumberCell 32 2 29 0.536111111111111
Set-TextCell 32 3 29 "None"
Set-TextCell 32 4 29 "Revise the assignment2"
Set-TextCell 32 5 29 "We found a new feature for assignment2. But there are still a lot of things left to finish."
Set-TextCell 32 6 29 "Since this week we talked to Kaj, he taught us a lot of things about report. Especially how to relate the diagrams to text, how to make report concise. "
Set-TextCell 32 7 29 "Not bad"
$ws.Rows.Item(32).RowHeight = 89

# ---- Row 33 ----
Set-TextCell 33 1 29 "2/17/2020"
Set-NumberCell 33 2 29 0.895138888888889
Set-TextCell 33 3 29 "None"
Set-TextCell 33 4 29 "Revise the assignment2"
Set-TextCell 33 5 29 "Add more stuff to the new feature and revise the feature2."
Set-TextCell 33 6 29 "When I look to our feature2 again, I find some places to improve. For example. I add a little more explanation to the diagram. What I can learn from this is that always consider yourself as a reader, could you totally understand this report?"
Set-TextCell 33 7 29 "Not bad"
$ws.Rows.Item(33).RowHeight = 123

# ---- Row 34 ----
Set-TextCell 34 1 29 "2/18/2020"
Set-NumberCell 34 2 29 0.966666666666667
Set-TextCell 34 3 29 "None"
Set-TextCell 34 4 29 "Start assignment3"
Set-TextCell 34 5 29 "We divide the assignment to three parts, each team member has one."
Set-TextCell 34 6 29 "I find out that I don't totally understand what's the stakeholders in our project. I still need more time on it. "
Set-TextCell 34 7 29 "Not bad"
$ws.Rows.Item(34).RowHeight = 85

[void]$ws.Range("B34").Select()
